$wb = $excel.ActiveWorkbook

# Sheets
$wsReadme = $wb.Worksheets.Item("README")
$wsBccw = $wb.Worksheets.Item("BME_BCCW")
$wsBcc = $wb.Worksheets.Item("BME_BCC")

# --- Update BME_BCCW (sheet2) data values ---
# Row 2
$wsBccw.Range("E2").Value = 173265.45
$wsBccw.Range("F2").Value = 366371.57
$wsBccw.Range("J2").Value = 47944
$wsBccw.Range("K2").Value = 264623

# Row 3
$wsBccw.Range("E3").Value = 167665.57
$wsBccw.Range("F3").Value = 457905.52
$wsBccw.Range("J3").Value = 47944
$wsBccw.Range("K3").Value = 264623

# Row 4
$wsBccw.Range("E4").Value = 4129.4799999999996
$wsBccw.Range("F4").Value = 365533.61
$wsBccw.Range("J4").Value = 47944
$wsBccw.Range("K4").Value = 264623

# Row 5
$wsBccw.Range("E5").Value = 1983.48
$wsBccw.Range("F5").Value = 385241.67
$wsBccw.Range("J5").Value = 272944
$wsBccw.Range("K5").Value = 389623

# Row 6
$wsBccw.Range("E6").Value = 4814.76
$wsBccw.Range("F6").Value = 544480.18000000005
$wsBccw.Range("J6").Value = 272944
$wsBccw.Range("K6").Value = 389623

# --- Update selections / views ---
$wsReadme.Select()
$wsReadme.Range("A4").Select()

$wsBccw.Select()
$wsBccw.Range("B5").Select()

$wsBcc.Select()
$wsBcc.Range("H22").Select()

# Make BME_BCC the active (selected) tab at the end
$wsBcc.Select()

$wb.Windows.Item(1).WindowState = -4143
